$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("B2").Value = 1.02
$ws.Range("C2").Value = 1.055467568591688
$ws.Range("D2").Value = 1.058801544484407
$ws.Range("E2").Value = 1.051880765516557
$ws.Range("F2").Value = 1.068206579710951
$ws.Range("I2").Value = 1.048956121735408
$ws.Range("J2").Value = 1.060474275598555
$ws.Range("K2").Value = 1.061532785797125
$ws.Range("L2").Value = 1.054631040917576
$ws.Range("M2").Value = 1.070912382260419
$ws.Range("B3").Value = 1.02
$ws.Range("C3").Value = 1.056599591360932
$ws.Range("D3").Value = 1.059682294013486
$ws.Range("E3").Value = 1.052851651317444
$ws.Range("F3").Value = 1.069228665748127
$ws.Range("I3").Value = 1.049292787907899
$ws.Range("J3").Value = 1.061257106485951
$ws.Range("K3").Value = 1.062227470949573
$ws.Range("L3").Value = 1.055414254920604
$ws.Range("M3").Value = 1.071749889252793
$ws.Range("B4").Value = 1.02
$ws.Range("C4").Value = 1.057332223478264
$ws.Range("D4").Value = 1.060252299544659
$ws.Range("E4").Value = 1.053480279486732
$ws.Range("F4").Value = 1.06989043968057
$ws.Range("I4").Value = 1.0495095323887
$ws.Range("J4").Value = 1.061763214695928
$ws.Range("K4").Value = 1.062676441140066
$ws.Range("L4").Value = 1.055920827520811
$ws.Range("M4").Value = 1.072291615656185
$ws.Range("B5").Value = 1.02
$ws.Range("C5").Value = 1.057640255550738
$ws.Range("D5").Value = 1.060491953959416
$ws.Range("E5").Value = 1.053744650411986
$ws.Range("F5").Value = 1.070168748844151
$ws.Range("I5").Value = 1.04960038814321
$ws.Range("J5").Value = 1.061975878321129
$ws.Range("K5").Value = 1.062865059446148
$ws.Range("L5").Value = 1.056133738110967
$ws.Range("M5").Value = 1.072519310263085
$ws.Range("B6").Value = 1.02
$ws.Range("C6").Value = 1.057691977473959
$ws.Range("D6").Value = 1.060532194395398
$ws.Range("E6").Value = 1.05378904505785
$ws.Range("F6").Value = 1.070215484001451
$ws.Range("I6").Value = 1.049615627764162
$ws.Range("J6").Value = 1.06201157936674
$ws.Range("K6").Value = 1.062896721737041
$ws.Range("L6").Value = 1.056169483654162
$ws.Range("M6").Value = 1.072557538420349
$ws.Range("B7").Value = 1.02
$ws.Range("C7").Value = 1.057336339285417
$ws.Range("D7").Value = 1.060255501724341
$ws.Range("E7").Value = 1.053483811648346
$ws.Range("F7").Value = 1.069894158070858
$ws.Range("I7").Value = 1.049510747443641
$ws.Range("J7").Value = 1.061766056727328
$ws.Range("K7").Value = 1.062678961973134
$ws.Range("L7").Value = 1.055923672649002
$ws.Range("M7").Value = 1.072294658308146
$ws.Range("B8").Value = 1.02
$ws.Range("C8").Value = 1.055850112817388
$ws.Range("D8").Value = 1.059099176701771
$ws.Range("E8").Value = 1.052208797977715
$ws.Range("F8").Value = 1.068551912250932
$ws.Range("I8").Value = 1.049070127603188
$ws.Range("J8").Value = 1.060738927208442
$ws.Range("K8").Value = 1.061767669158967
$ws.Range("L8").Value = 1.054895777122226
$ws.Range("M8").Value = 1.071195462320851
$ws.Range("B9").Value = 1.02
$ws.Range("C9").Value = 1.053232220787134
$ws.Range("D9").Value = 1.057062373996698
$ws.Range("E9").Value = 1.049965131740276
$ws.Range("F9").Value = 1.066189896155058
$ws.Range("I9").Value = 1.048285267724917
$ws.Range("J9").Value = 1.058925654216182
$ws.Range("K9").Value = 1.06015774049203
$ws.Range("L9").Value = 1.053082816477076
$ws.Range("M9").Value = 1.069257036717271
$ws.Range("B10").Value = 1.02
$ws.Range("C10").Value = 1.051487612693381
$ws.Range("D10").Value = 1.055705048729484
$ws.Range("E10").Value = 1.048471424522321
$ws.Range("F10").Value = 1.064617371978243
$ws.Range("I10").Value = 1.047756359348819
$ws.Range("J10").Value = 1.057714554322024
$ws.Range("K10").Value = 1.059081689918882
$ws.Range("L10").Value = 1.051873048047128
$ws.Range("M10").Value = 1.067963746180545
$ws.Range("B11").Value = 1.02
$ws.Range("C11").Value = 1.050732323654786
$ws.Range("D11").Value = 1.055117442527879
$ws.Range("E11").Value = 1.047825124168228
$ws.Range("F11").Value = 1.063936962874926
$ws.Range("I11").Value = 1.047525991316724
$ws.Range("J11").Value = 1.057189598627653
$ws.Range("K11").Value = 1.058615091920446
$ws.Range("L11").Value = 1.051348936564865
$ws.Range("M11").Value = 1.067403498088884
$ws.Range("B12").Value = 1.02
$ws.Range("C12").Value = 1.050451795459563
$ws.Range("D12").Value = 1.054899198188402
$ws.Range("E12").Value = 1.047585132188328
$ws.Range("F12").Value = 1.063684304329437
$ws.Range("I12").Value = 1.047440219868685
$ws.Range("J12").Value = 1.05699452508271
$ws.Range("K12").Value = 1.058441677207331
$ws.Range("L12").Value = 1.051154216924433
$ws.Range("M12").Value = 1.067195360109713
$ws.Range("B13").Value = 1.02
$ws.Range("C13").Value = 1.050511968829176
$ws.Range("D13").Value = 1.054946011490057
$ws.Range("E13").Value = 1.047636607991912
$ws.Range("F13").Value = 1.063738497021004
$ws.Range("I13").Value = 1.047458627300755
$ws.Range("J13").Value = 1.057036372718538
$ws.Range("K13").Value = 1.058478879754022
$ws.Range("L13").Value = 1.051195986819909
$ws.Range("M13").Value = 1.067240008084203
$ws.Range("B14").Value = 1.02
$ws.Range("C14").Value = 1.05070913470838
$ws.Range("D14").Value = 1.055099401999581
$ws.Range("E14").Value = 1.047805284872707
$ws.Range("F14").Value = 1.063916076497051
$ws.Range("I14").Value = 1.04751890555366
$ws.Range("J14").Value = 1.057173475454342
$ws.Range("K14").Value = 1.058600759435502
$ws.Range("L14").Value = 1.051332841822124
$ws.Range("M14").Value = 1.067386294091898
$ws.Range("B15").Value = 1.02
$ws.Range("C15").Value = 1.050830617707934
$ws.Range("D15").Value = 1.055193913443263
$ws.Range("E15").Value = 1.047909221903554
$ws.Range("F15").Value = 1.064025499100902
$ws.Range("I15").Value = 1.047556018138039
$ws.Range("J15").Value = 1.057257938150201
$ws.Range("K15").Value = 1.058675840348902
$ws.Range("L15").Value = 1.051417157228103
$ws.Range("M15").Value = 1.06747642083702
$ws.Range("B16").Value = 1.02
$ws.Range("C16").Value = 1.051537741563889
$ws.Range("D16").Value = 1.055744048847024
$ws.Range("E16").Value = 1.048514327554379
$ws.Range("F16").Value = 1.064662539107706
$ws.Range("I16").Value = 1.047771619723558
$ws.Range("J16").Value = 1.057749382470934
$ws.Range("K16").Value = 1.059112642553409
$ws.Range("L16").Value = 1.051907825817623
$ws.Range("M16").Value = 1.068000922834174
$ws.Range("B17").Value = 1.02
$ws.Range("C17").Value = 1.051981337916041
$ws.Range("D17").Value = 1.05608916764325
$ws.Range("E17").Value = 1.048894024063793
$ws.Range("F17").Value = 1.065062272426919
$ws.Range("I17").Value = 1.047906500143221
$ws.Range("J17").Value = 1.058057507350566
$ws.Range("K17").Value = 1.059386459973351
$ws.Range("L17").Value = 1.052215535845268
$ws.Range("M17").Value = 1.06832986346624
$ws.Range("B18").Value = 1.02
$ws.Range("C18").Value = 1.052240093518146
$ws.Range("D18").Value = 1.056290481616669
$ws.Range("E18").Value = 1.049115541388169
$ws.Range("F18").Value = 1.06529547876554
$ws.Range("I18").Value = 1.047985043629899
$ws.Range("J18").Value = 1.058237179156595
$ws.Range("K18").Value = 1.059546109273819
$ws.Range("L18").Value = 1.05239499149493
$ws.Range("M18").Value = 1.068521705364476
$ws.Range("B19").Value = 1.02
$ws.Range("C19").Value = 1.052328324784787
$ws.Range("D19").Value = 1.056359126516277
$ws.Range("E19").Value = 1.049191081007348
$ws.Range("F19").Value = 1.065375004355414
$ws.Range("I19").Value = 1.048011802885299
$ws.Range("J19").Value = 1.058298433749893
$ws.Range("K19").Value = 1.059600534747497
$ws.Range("L19").Value = 1.052456176760597
$ws.Range("M19").Value = 1.068587114493645
$ws.Range("B20").Value = 1.02
$ws.Range("C20").Value = 1.051933742850694
$ws.Range("D20").Value = 1.056052138419005
$ws.Range("E20").Value = 1.04885328138409
$ws.Range("F20").Value = 1.065019379792397
$ws.Range("I20").Value = 1.047892042195805
$ws.Range("J20").Value = 1.058024453870645
$ws.Range("K20").Value = 1.059357088552259
$ws.Range("L20").Value = 1.052182524203586
$ws.Range("M20").Value = 1.068294573692054
$ws.Range("B21").Value = 1.02
$ws.Range("C21").Value = 1.05065107373242
$ws.Range("D21").Value = 1.055054231799238
$ws.Range("E21").Value = 1.047755611727881
$ws.Range("F21").Value = 1.063863781678497
$ws.Range("I21").Value = 1.047501160701462
$ws.Range("J21").Value = 1.05713310437137
$ws.Range("K21").Value = 1.058564871653206
$ws.Range("L21").Value = 1.051292542576423
$ws.Range("M21").Value = 1.067343217528104
$ws.Range("B22").Value = 1.02
$ws.Range("C22").Value = 1.049844722755846
$ws.Range("D22").Value = 1.054426916904259
$ws.Range("E22").Value = 1.047065884109382
$ws.Range("F22").Value = 1.063137649169858
$ws.Range("I22").Value = 1.047254226319574
$ws.Range("J22").Value = 1.056572205161163
$ws.Range("K22").Value = 1.058066198044439
$ws.Range("L22").Value = 1.050732736773056
$ws.Range("M22").Value = 1.066744848638752
$ws.Range("B23").Value = 1.02
$ws.Range("C23").Value = 1.050272174163858
$ws.Range("D23").Value = 1.054759458088961
$ws.Range("E23").Value = 1.047431481902833
$ws.Range("F23").Value = 1.063522544266394
$ws.Range("I23").Value = 1.047385242001423
$ws.Range("J23").Value = 1.056869593316434
$ws.Range("K23").Value = 1.058330608945361
$ws.Range("L23").Value = 1.051029523134932
$ws.Range("M23").Value = 1.067062075586979
$ws.Range("B24").Value = 1.02
$ws.Range("C24").Value = 1.05195524897044
$ws.Range("D24").Value = 1.056068870294288
$ws.Range("E24").Value = 1.048871691102633
$ws.Range("F24").Value = 1.065038760977161
$ws.Range("I24").Value = 1.047898575521437
$ws.Range("J24").Value = 1.058039389476641
$ws.Range("K24").Value = 1.059370360430069
$ws.Range("L24").Value = 1.052197440824928
$ws.Range("M24").Value = 1.068310519694776
$ws.Range("B25").Value = 1.02
$ws.Range("C25").Value = 1.053908890676412
$ws.Range("D25").Value = 1.057588841468961
$ws.Range("E25").Value = 1.050544808758574
$ws.Range("F25").Value = 1.066800154908182
$ws.Range("I25").Value = 1.048489171952858
$ws.Range("J25").Value = 1.059394825114732
$ws.Range("K25").Value = 1.060574433063002
$ws.Range("L25").Value = 1.053551709064183
$ws.Range("M25").Value = 1.069758343873268
